# Adds the 2022 (column S) data point to the "Suicide mortality rate" table.
# Mirrors: +1 new year column appended after R (2021), for rows 4 (year
# header) through 14 (last data row / "г.Ош"), plus updated selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 4: year header "2022" — same look as the other year-header cells
# (bold Times New Roman 9pt, centred vertically / right horizontally,
# medium border above & below, matching R4's "2021").
# ---------------------------------------------------------------------
$c = $ws.Range("S4")
$c.Value = 2022
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.Font.Bold = $true
$c.Font.Family = 1
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Borders.Item(8).LineStyle = -4138
$c.Borders.Item(8).Weight = -4138
$c.Borders.Item(9).LineStyle = -4138
$c.Borders.Item(9).Weight = -4138

# ---------------------------------------------------------------------
# Row 5: totals row ("Кыргызская Республика") — bold font, "0.0" format
# ---------------------------------------------------------------------
$c = $ws.Range("S5")
$c.Value = 4.9000000000000004
$c.NumberFormat = "0.0"
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.Font.Bold = $true
$c.Font.Family = 1
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Rows 6-13: regional data rows — regular (non-bold) font, "0.0" format
# ---------------------------------------------------------------------
$values = @{
    6  = 3.4
    7  = 3.5
    8  = 13.1
    9  = 8.1
    10 = 2.5
    11 = 2.6
    12 = 10.8
    13 = 2.1
}

foreach ($row in 6..13) {
    $c = $ws.Range("S$row")
    $c.Value = $values[$row]
    $c.NumberFormat = "0.0"
    $c.Font.Name = "Times New Roman"
    $c.Font.Size = 9
    $c.Font.Bold = $false
    $c.Font.Family = 1
    $c.HorizontalAlignment = -4152
    $c.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------
# Row 14: last data row ("г.Ош") — same as rows 6-13 but with a medium
# bottom border closing off the table (matches R14).
# ---------------------------------------------------------------------
$c = $ws.Range("S14")
$c.Value = 1.1000000000000001
$c.NumberFormat = "0.0"
$c.Font.Name = "Times New Roman"
$c.Font.Size = 9
$c.Font.Bold = $false
$c.Font.Family = 1
$c.HorizontalAlignment = -4152
$c.VerticalAlignment = -4108
$c.Borders.Item(9).LineStyle = -4138
$c.Borders.Item(9).Weight = -4138

# ---------------------------------------------------------------------
# Selection moves from T9 to T4 (right next to the new data column).
# ---------------------------------------------------------------------
$ws.Range("T4").Select()
